# Automatische test-sync: 2025-06-19 21:39:50
# Appends a new "Factuur verzoek" log entry to the Logs sheet and
# refreshes the Dashboard count for "Factuur / Administratie".

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$newRow = 22

$logs.Cells.Item($newRow, 1).Value = "Factuur verzoek"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Kunt u mij de factuur van mijn laatste bestelling toesturen?"
$logs.Cells.Item($newRow, 4).Value = "Factuur / Administratie"
$logs.Cells.Item($newRow, 6).Value = "2025-06-19 21:39:10"
$logs.Cells.Item($newRow, 7).Value = "Nee"

# Extend the conditional-formatting ranges so the new row is covered too
# (Categorie column D and Beantwoord column G), keeping every existing
# rule's priority/format untouched.
$catRange = $logs.Range("D2:D" + $newRow)
foreach ($fc in $catRange.FormatConditions) {
    $fc.ModifyAppliesToRange($catRange)
}

$answeredRange = $logs.Range("G2:G" + $newRow)
foreach ($fc in $answeredRange.FormatConditions) {
    $fc.ModifyAppliesToRange($answeredRange)
}

# Bump the Dashboard tally for "Factuur / Administratie" from 2 to 3.
$dashboard.Cells.Item(4, 2).Value = 3
